# Refresh the "NDC No Hyphens" FHIR ValueSet metadata sheet for the 6.0.0
# release (Alvearie IG build from 2022-01-21).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: Property / Value table ----------------------------------------
# (Sheet2 "Codes" table is unaffected by this edit; its shared-string index
# shift is handled automatically when the strings above are rewritten.)

# Version bump: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Refresh the IG build date/time
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Rows 11-15 (the second duplicate "Contact" row through "Immutable") shift
# up one row to 10-14, collapsing the accidental duplicate Contact row and
# shrinking the table from 15 to 14 rows.
$shiftedValues = $ws1.Range("A11:B15").Value2
$ws1.Range("A10:B14").Value = $shiftedValues
$ws1.Rows.Item(15).Delete()

# Publisher now has a value, and the old duplicate "Contact" row (now row 10,
# after the shift above) becomes the new "Jurisdiction" row.
$ws1.Range("B9").Value = "Alvearie Team"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"
